$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "image"
$ws.Cells.Item(1, 2).Value = "trialnum"
$ws.Cells.Item(1, 3).Value = "condition"
$ws.Cells.Item(1, 4).Value = "word"
$ws.Cells.Item(1, 5).Value = "location"
$ws.Cells.Item(1, 6).Value = "repetition"

# Update existing data rows 2-25 and append new data rows 26-49
$ws.Cells.Item(2, 1).Value = "D.png"
$ws.Cells.Item(2, 2).Value = 241
$ws.Cells.Item(2, 3).Value = "R"
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 1).Value = "K.png"
$ws.Cells.Item(3, 2).Value = 242
$ws.Cells.Item(3, 3).Value = "R"
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 1).Value = "K.png"
$ws.Cells.Item(4, 2).Value = 243
$ws.Cells.Item(4, 3).Value = "R"
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(5, 1).Value = "H.png"
$ws.Cells.Item(5, 2).Value = 244
$ws.Cells.Item(5, 3).Value = "R"
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 1).Value = "H.png"
$ws.Cells.Item(6, 2).Value = 245
$ws.Cells.Item(6, 3).Value = "R"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(7, 1).Value = "J.png"
$ws.Cells.Item(7, 2).Value = 246
$ws.Cells.Item(7, 3).Value = "R"
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(8, 1).Value = "F.png"
$ws.Cells.Item(8, 2).Value = 247
$ws.Cells.Item(8, 3).Value = "R"
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 1).Value = "M.png"
$ws.Cells.Item(9, 2).Value = 248
$ws.Cells.Item(9, 3).Value = "R"
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 1).Value = "B.png"
$ws.Cells.Item(10, 2).Value = 249
$ws.Cells.Item(10, 3).Value = "R"
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(11, 1).Value = "C.png"
$ws.Cells.Item(11, 2).Value = 250
$ws.Cells.Item(11, 3).Value = "R"
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(12, 1).Value = "C.png"
$ws.Cells.Item(12, 2).Value = 251
$ws.Cells.Item(12, 3).Value = "R"
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(13, 1).Value = "E.png"
$ws.Cells.Item(13, 2).Value = 252
$ws.Cells.Item(13, 3).Value = "R"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(14, 1).Value = "L.png"
$ws.Cells.Item(14, 2).Value = 253
$ws.Cells.Item(14, 3).Value = "R"
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(15, 1).Value = "G.png"
$ws.Cells.Item(15, 2).Value = 254
$ws.Cells.Item(15, 3).Value = "R"
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(16, 1).Value = "M.png"
$ws.Cells.Item(16, 2).Value = 255
$ws.Cells.Item(16, 3).Value = "R"
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(17, 1).Value = "B.png"
$ws.Cells.Item(17, 2).Value = 256
$ws.Cells.Item(17, 3).Value = "R"
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(18, 1).Value = "E.png"
$ws.Cells.Item(18, 2).Value = 257
$ws.Cells.Item(18, 3).Value = "R"
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(19, 1).Value = "L.png"
$ws.Cells.Item(19, 2).Value = 258
$ws.Cells.Item(19, 3).Value = "R"
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(20, 1).Value = "E.png"
$ws.Cells.Item(20, 2).Value = 259
$ws.Cells.Item(20, 3).Value = "R"
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(21, 1).Value = "H.png"
$ws.Cells.Item(21, 2).Value = 260
$ws.Cells.Item(21, 3).Value = "R"
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(22, 1).Value = "D.png"
$ws.Cells.Item(22, 2).Value = 261
$ws.Cells.Item(22, 3).Value = "R"
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(23, 1).Value = "A.png"
$ws.Cells.Item(23, 2).Value = 262
$ws.Cells.Item(23, 3).Value = "R"
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(24, 1).Value = "G.png"
$ws.Cells.Item(24, 2).Value = 263
$ws.Cells.Item(24, 3).Value = "R"
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(25, 1).Value = "L.png"
$ws.Cells.Item(25, 2).Value = 264
$ws.Cells.Item(25, 3).Value = "R"
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(26, 1).Value = "D.png"
$ws.Cells.Item(26, 2).Value = 265
$ws.Cells.Item(26, 3).Value = "R"
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(27, 1).Value = "A.png"
$ws.Cells.Item(27, 2).Value = 266
$ws.Cells.Item(27, 3).Value = "R"
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 2
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(28, 1).Value = "H.png"
$ws.Cells.Item(28, 2).Value = 267
$ws.Cells.Item(28, 3).Value = "R"
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(29, 1).Value = "G.png"
$ws.Cells.Item(29, 2).Value = 268
$ws.Cells.Item(29, 3).Value = "R"
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(30, 1).Value = "G.png"
$ws.Cells.Item(30, 2).Value = 269
$ws.Cells.Item(30, 3).Value = "R"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(31, 1).Value = "M.png"
$ws.Cells.Item(31, 2).Value = 270
$ws.Cells.Item(31, 3).Value = "R"
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(32, 1).Value = "F.png"
$ws.Cells.Item(32, 2).Value = 271
$ws.Cells.Item(32, 3).Value = "R"
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(33, 1).Value = "M.png"
$ws.Cells.Item(33, 2).Value = 272
$ws.Cells.Item(33, 3).Value = "R"
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(34, 1).Value = "K.png"
$ws.Cells.Item(34, 2).Value = 273
$ws.Cells.Item(34, 3).Value = "R"
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(35, 1).Value = "K.png"
$ws.Cells.Item(35, 2).Value = 274
$ws.Cells.Item(35, 3).Value = "R"
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 1
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(36, 1).Value = "J.png"
$ws.Cells.Item(36, 2).Value = 275
$ws.Cells.Item(36, 3).Value = "R"
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(37, 1).Value = "B.png"
$ws.Cells.Item(37, 2).Value = 276
$ws.Cells.Item(37, 3).Value = "R"
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(38, 1).Value = "B.png"
$ws.Cells.Item(38, 2).Value = 277
$ws.Cells.Item(38, 3).Value = "R"
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(39, 1).Value = "C.png"
$ws.Cells.Item(39, 2).Value = 278
$ws.Cells.Item(39, 3).Value = "R"
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 2
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(40, 1).Value = "F.png"
$ws.Cells.Item(40, 2).Value = 279
$ws.Cells.Item(40, 3).Value = "R"
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 3
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(41, 1).Value = "F.png"
$ws.Cells.Item(41, 2).Value = 280
$ws.Cells.Item(41, 3).Value = "R"
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 1
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(42, 1).Value = "J.png"
$ws.Cells.Item(42, 2).Value = 281
$ws.Cells.Item(42, 3).Value = "R"
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 2
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(43, 1).Value = "A.png"
$ws.Cells.Item(43, 2).Value = 282
$ws.Cells.Item(43, 3).Value = "R"
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 3
$ws.Cells.Item(43, 6).Value = 1
$ws.Cells.Item(44, 1).Value = "C.png"
$ws.Cells.Item(44, 2).Value = 283
$ws.Cells.Item(44, 3).Value = "R"
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 1
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(45, 1).Value = "A.png"
$ws.Cells.Item(45, 2).Value = 284
$ws.Cells.Item(45, 3).Value = "R"
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 2
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(46, 1).Value = "L.png"
$ws.Cells.Item(46, 2).Value = 285
$ws.Cells.Item(46, 3).Value = "R"
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 3
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(47, 1).Value = "J.png"
$ws.Cells.Item(47, 2).Value = 286
$ws.Cells.Item(47, 3).Value = "R"
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 1
$ws.Cells.Item(47, 6).Value = 1
$ws.Cells.Item(48, 1).Value = "D.png"
$ws.Cells.Item(48, 2).Value = 287
$ws.Cells.Item(48, 3).Value = "R"
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 2
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(49, 1).Value = "E.png"
$ws.Cells.Item(49, 2).Value = 288
$ws.Cells.Item(49, 3).Value = "R"
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 3
$ws.Cells.Item(49, 6).Value = 1

# Update selection to match the new active range
$ws.Range("A26:F49").Select() | Out-Null
